$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update employee numbers for Bravo, Charlie, Delta
$ws.Range("B3").Value = "B02"
$ws.Range("B4").Value = "C03"
$ws.Range("B5").Value = "D04"

# Replace "Freddie" row (row 6) with "Echo" row
$ws.Range("B6").Value = "E05"
$ws.Range("C6").Value = "echo@msn.com"
$ws.Range("A6").Value = "Echo"
